# Add team record (Wins/Losses/Ties) columns to the roster sheet.
# The W/L/T values are the same for every player row (team-level record
# repeated down the column), per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled like the rest of row 1 (bold/centered/bordered).
# Copy the formatting from the last existing header cell (AB1) so the new
# headers reuse the same cell style instead of minting a near-duplicate one.
$ws.Range("AB1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record repeated for every data row (2-45): 77 wins, 85 losses, 0 ties.
$ws.Range("AD2:AD45").Value = 77
$ws.Range("AE2:AE45").Value = 85
$ws.Range("AF2:AF45").Value = 0
